$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column G (header 'K') values per row, per the regenerated save_data
$gValues = @{
    2 = 4
    3 = 1
    4 = 1
    5 = 1
    6 = 3
    7 = 3
    8 = 0
    9 = 0
    10 = 0
    11 = 0
    12 = 1
    13 = 2
    14 = 2
    15 = 3
    16 = 2
    17 = 0
    18 = 1
    19 = 2
    20 = 1
    21 = 1
    22 = 1
    23 = 2
    24 = 1
    25 = 1
    26 = 1
    27 = 2
    28 = 1
    29 = 0
    30 = 1
    31 = 2
    32 = 1
    33 = 1
    34 = 1
    35 = 2
    36 = 1
    37 = 1
    38 = 1
    39 = 0
    40 = 2
    41 = 1
    42 = 0
    43 = 1
    44 = 1
    45 = 0
    46 = 3
    47 = 1
    48 = 3
    49 = 0
    50 = 1
    51 = 1
    52 = 1
    53 = 1
    54 = 1
    55 = 1
    56 = 1
    57 = 1
    58 = 3
    59 = 0
    60 = 0
    61 = 2
    62 = 1
    63 = 1
    64 = 0
    65 = 2
    66 = 1
    67 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $gValues[$row]
}

Write-Output "Updated $($gValues.Count) cells in column G"
